# "Reference to dxfViewer added and Status updated"
#
# The "Mechanical" row (row 18) in Sheet1's Status column previously read
# "Needs reference to dxf file reader". A reference to the dxf viewer has
# now been added, so the status is updated to "Complete". Since that was
# the only cell using that particular shared string, the string itself
# drops out of the shared-strings table on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C18").Value = "Complete"

# Reflect the author's cursor / scroll position at the time of the edit:
# selection moved down to C19, with the view scrolled so row 3 is at the top.
$ws.Range("C19").Select()
$excel.ActiveWindow.ScrollRow = 3
